$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 219: fill in newly-populated numeric cells (columns D,E,M,O,P,Q) ---
$ws.Range("D219").Value = -0.1
$ws.Range("E219").Value = -0.5
$ws.Range("M219").Value = 1.13
$ws.Range("O219").Value = 38
$ws.Range("P219").Value = 7.75
$ws.Range("Q219").Value = 2.5

# --- New rows 220-223: dates go to column A (text, not Excel dates) ---
$ws.Range("A220:A223").NumberFormat = "@"
$ws.Range("A220").Value = "03-11-2021"
$ws.Range("A221").Value = "04-11-2021"
$ws.Range("A222").Value = "05-11-2021"
$ws.Range("A223").Value = "08-11-2021"
$ws.Range("A220:A223").ClearFormats()

# --- Row 220 ---
$ws.Range("B220").Value = 0.25
$ws.Range("C220").Value = 0.1
$ws.Range("D220").Value = -0.1
$ws.Range("E220").Value = -0.5
$ws.Range("F220").Value = 0.75
$ws.Range("G220").Value = 4.35
$ws.Range("H220").Value = 1.75
$ws.Range("I220").Value = 0.5
$ws.Range("J220").Value = 1.5
$ws.Range("K220").Value = 7.5
$ws.Range("L220").Value = 0.5
$ws.Range("M220").Value = 1.13
$ws.Range("N220").Value = 16
$ws.Range("O220").Value = 38
$ws.Range("P220").Value = 7.75
$ws.Range("Q220").Value = 2.5
$ws.Range("R220").Value = 4.75
$ws.Range("S220").Value = 1.5

# --- Row 221 ---
$ws.Range("B221").Value = 0.25
$ws.Range("C221").Value = 0.1
$ws.Range("D221").Value = -0.1
$ws.Range("E221").Value = -0.5
$ws.Range("F221").Value = 0.75
$ws.Range("G221").Value = 4.35
$ws.Range("H221").Value = 1.75
$ws.Range("I221").Value = 1.25
$ws.Range("J221").Value = 1.5
$ws.Range("L221").Value = 0.5
$ws.Range("M221").Value = 1.13
$ws.Range("N221").Value = 16
$ws.Range("O221").Value = 38
$ws.Range("P221").Value = 7.75
$ws.Range("Q221").Value = 2.5
$ws.Range("R221").Value = 4.75
$ws.Range("S221").Value = 1.5

# --- Row 222 ---
$ws.Range("B222").Value = 0.25
$ws.Range("C222").Value = 0.1
$ws.Range("D222").Value = -0.1
$ws.Range("E222").Value = -0.5
$ws.Range("F222").Value = 0.75
$ws.Range("G222").Value = 4.35
$ws.Range("H222").Value = 1.75
$ws.Range("I222").Value = 1.25
$ws.Range("J222").Value = 2.75
$ws.Range("K222").Value = 7.5
$ws.Range("L222").Value = 0.5
$ws.Range("M222").Value = 1.13
$ws.Range("N222").Value = 16
$ws.Range("O222").Value = 38
$ws.Range("P222").Value = 7.75
$ws.Range("Q222").Value = 2.5
$ws.Range("R222").Value = 4.75
$ws.Range("S222").Value = 1.5

# --- Row 223 ---
$ws.Range("C223").Value = 0.1
$ws.Range("F223").Value = 0.75
$ws.Range("G223").Value = 4.35
$ws.Range("I223").Value = 1.25
$ws.Range("J223").Value = 2.75
$ws.Range("K223").Value = 7.5
$ws.Range("L223").Value = 0.5
$ws.Range("N223").Value = 16
$ws.Range("Q223").Value = 2.5
$ws.Range("R223").Value = 4.75
$ws.Range("S223").Value = 1.5
